$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ligand/receptor TPM-derived values for the Icam1-Itgax sheet.
# Only the cells listed below change; all other data is untouched.

# Row 2
$ws.Range("G2").Value = 9.861094666666666
$ws.Range("H2").Value = 29.583284
$ws.Range("I2").Value = 0.243709096397741
$ws.Range("J2").Value = 0.2437090963977409
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.007501333333333333
$ws.Range("N2").Value = 0.022504
$ws.Range("O2").Value = 0.1758070060310615
$ws.Range("P2").Value = 0.1758070060310615
$ws.Range("Q2").Value = 0.07397135812622221
$ws.Range("R2").Value = 0.665742223136
$ws.Range("S2").Value = 0.0428457665802222
$ws.Range("T2").Value = 0.0428457665802222

# Row 3
$ws.Range("G3").Value = 9.861094666666666
$ws.Range("H3").Value = 29.583284
$ws.Range("I3").Value = 0.243709096397741
$ws.Range("J3").Value = 0.2437090963977409
$ws.Range("O3").Value = 0.8241929939689385
$ws.Range("P3").Value = 0.8241929939689384
$ws.Range("Q3").Value = 0.3467818291111111
$ws.Range("R3").Value = 3.121036462
$ws.Range("S3").Value = 0.2008633298175188
$ws.Range("T3").Value = 0.2008633298175187

# Row 4
$ws.Range("I4").Value = 0.7254466225154019
$ws.Range("J4").Value = 0.7254466225154018
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.007501333333333333
$ws.Range("N4").Value = 0.022504
$ws.Range("O4").Value = 0.1758070060310615
$ws.Range("P4").Value = 0.1758070060310615
$ws.Range("Q4").Value = 0.2201898604062222
$ws.Range("R4").Value = 1.981708743656
$ws.Range("S4").Value = 0.1275385987397785
$ws.Range("T4").Value = 0.1275385987397785

# Row 5
$ws.Range("I5").Value = 0.7254466225154019
$ws.Range("J5").Value = 0.7254466225154018
$ws.Range("O5").Value = 0.8241929939689385
$ws.Range("P5").Value = 0.8241929939689384
$ws.Range("S5").Value = 0.5979080237756235
$ws.Range("T5").Value = 0.5979080237756234

# Row 6
$ws.Range("I6").Value = 0.03084428108685718
$ws.Range("J6").Value = 0.03084428108685716
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.007501333333333333
$ws.Range("N6").Value = 0.022504
$ws.Range("O6").Value = 0.1758070060310615
$ws.Range("P6").Value = 0.1758070060310615
$ws.Range("Q6").Value = 0.009361954051555555
$ws.Range("R6").Value = 0.084257586464
$ws.Range("S6").Value = 0.005422640711060856
$ws.Range("T6").Value = 0.005422640711060855

# Row 7
$ws.Range("I7").Value = 0.03084428108685718
$ws.Range("J7").Value = 0.03084428108685716
$ws.Range("O7").Value = 0.8241929939689385
$ws.Range("P7").Value = 0.8241929939689384
$ws.Range("S7").Value = 0.02542164037579632
$ws.Range("T7").Value = 0.02542164037579631
